# Total number of domain's count
# Adds a grand-total row (SUM) under the existing data, widens the two
# data columns so the longer domain names/labels are fully visible, and
# moves the selection to the newly added total cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New total row: sum of all the per-domain configuration counts in column B.
$ws.Range("B68").Formula = "=SUM(B2:B67)"

# Widen columns A and B to comfortably fit the (now longer-looking) content.
$ws.Columns("A").ColumnWidth = 32.666666666666664
$ws.Columns("B").ColumnWidth = 32.5

# Move/extend the active selection down to the new total cell, as the
# author would have after typing the formula.
[void]$ws.Range("B68").Select()
